# Update cryptocurrency price/volume data per Oct 21 2023 GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number need to be forced to
# stay as text (matching the original inlineStr cells) instead of being
# auto-converted to a numeric value by Excel.
$textCells = @(
    'D5',
    'D6',
    'D7',
    'D8',
    'D9',
    'D10',
    'D14',
    'D18',
    'D19',
    'D21',
    'D23',
    'D25',
    'D26',
    'D32',
    'D33',
    'D37',
    'D38',
    'D39',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D50',
    'D51'
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.725.23'
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').Value = '1.619.79'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').Value = '  -0.61%  '
$ws.Range('D5').Value = '212.79'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').Value = '0.522'
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('D7').Value = '0.990'
$ws.Range('E7').Value = '  -0.69%  '
$ws.Range('D8').Value = '29.25'
$ws.Range('E8').Value = '  +9.84%  '
$ws.Range('D9').Value = '0.259'
$ws.Range('E9').Value = '  +3.31%  '
$ws.Range('D10').Value = '0.0609'
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('D12').Value = '1.847.11'
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('D13').Value = '1.600.96'
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('D14').Value = '0.569'
$ws.Range('E14').Value = '  +6.77%  '
$ws.Range('E15').Value = '  +5.86%  '
$ws.Range('E16').Value = '  +18.20%  '
$ws.Range('D17').Value = '29.712.17'
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('D18').Value = '64.26'
$ws.Range('E18').Value = '  +1.76%  '
$ws.Range('D19').Value = '242.06'
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('E20').Value = '  +3.34%  '
$ws.Range('D21').Value = '0.993'
$ws.Range('E21').Value = '  -0.46%  '
$ws.Range('E22').Value = '  +3.01%  '
$ws.Range('D23').Value = '9.71'
$ws.Range('E23').Value = '  +6.23%  '
$ws.Range('E24').Value = '  +1.29%  '
$ws.Range('D25').Value = '156.37'
$ws.Range('E25').Value = '  +1.30%  '
$ws.Range('D26').Value = '15.69'
$ws.Range('E26').Value = '  +2.78%  '
$ws.Range('E27').Value = '  +2.00%  '
$ws.Range('E28').Value = '  +3.63%  '
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('E30').Value = '  +3.46%  '
$ws.Range('E31').Value = '  +2.85%  '
$ws.Range('D32').Value = '3.32'
$ws.Range('E32').Value = '  +3.34%  '
$ws.Range('D33').Value = '3.21'
$ws.Range('E33').Value = '  +3.88%  '
$ws.Range('D34').Value = '1.425.68'
$ws.Range('E35').Value = '  +6.86%  '
$ws.Range('E36').Value = '  +1.36%  '
$ws.Range('D37').Value = '2.88'
$ws.Range('E37').Value = '  +2.47%  '
$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').Value = '2.28'
$ws.Range('E38').Value = '  -0.76%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.0171'
$ws.Range('E39').Value = '  +3.29%  '
$ws.Range('E40').Value = '  +4.27%  '
$ws.Range('D41').Value = '0.0505'
$ws.Range('E41').Value = '  +2.87%  '
$ws.Range('D42').Value = '0.831'
$ws.Range('E42').Value = '  +4.41%  '
$ws.Range('D43').Value = '1.97'
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').Value = '70.01'
$ws.Range('E44').Value = '  +6.73%  '
$ws.Range('D45').Value = '53.70'
$ws.Range('E45').Value = '  +1.90%  '
$ws.Range('D46').Value = '0.991'
$ws.Range('E46').Value = '  -0.67%  '
$ws.Range('E47').Value = '  +18.70%  '
$ws.Range('E48').Value = '  +3.60%  '
$ws.Range('D49').Value = '1.756.81'
$ws.Range('E49').Value = '  +0.63%  '
$ws.Range('D50').Value = '88.03'
$ws.Range('E50').Value = '  +1.64%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.0534'
$ws.Range('E51').Value = '  +2.30%  '

# Restore default (unstyled) appearance on the cells we forced to text format
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
